# Update scripts with new TPM values.
# - Sending cluster is now only "FAPs" (rows for "ECs" sending cluster are removed).
# - The three target-cluster rows (ECs, FAPs, MuSCs) survive with recomputed values.
# - Data now occupies rows 2-4 instead of 2-7; rows 5-7 are deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Dll3/Notch3 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Dll3"
$ws.Range("C2").Value = "Notch3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2138853333333333
$ws.Range("H2").Value = 0.641656
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.413580666666667
$ws.Range("N2").Value = 22.240742
$ws.Range("O2").Value = 0.05108888817597561
$ws.Range("P2").Value = 0.05108888817597561
$ws.Range("Q2").Value = 1.585656172083556
$ws.Range("R2").Value = 14.270905548752
$ws.Range("S2").Value = 0.05108888817597561
$ws.Range("T2").Value = 0.05108888817597561

# Row 3: FAPs -> Dll3/Notch3 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Dll3"
$ws.Range("C3").Value = "Notch3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2138853333333333
$ws.Range("H3").Value = 0.641656
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.0112127179963522
$ws.Range("P3").Value = 0.0112127179963522
$ws.Range("Q3").Value = 0.3480113999644445
$ws.Range("R3").Value = 3.13210259968
$ws.Range("S3").Value = 0.0112127179963522
$ws.Range("T3").Value = 0.0112127179963522

# Row 4: FAPs -> Dll3/Notch3 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Dll3"
$ws.Range("C4").Value = "Notch3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2138853333333333
$ws.Range("H4").Value = 0.641656
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9376983938276722
$ws.Range("P4").Value = 0.9376983938276722
$ws.Range("Q4").Value = 29.10353501145245
$ws.Range("R4").Value = 261.931815103072
$ws.Range("S4").Value = 0.9376983938276722
$ws.Range("T4").Value = 0.9376983938276722

# Rows 5-7 (old "ECs" sending-cluster data) no longer exist - delete them
# so remaining rows shift up and the used range shrinks to A1:T4.
$ws.Range("A5:T7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
